$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text-typed like the original inline strings
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.679.81'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '1.605.56'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '212.43'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '28.93'
$ws.Range("E8").Value = '  +7.68%  '
$ws.Range("E9").Value = '  +3.02%  '
$ws.Range("D10").Value = '0.0608'
$ws.Range("E10").Value = '  +1.97%  '
$ws.Range("D11").Value = '0.0907'
$ws.Range("D12").Value = '1.834.24'
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").Value = '1.610.22'
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").Value = '0.559'
$ws.Range("E14").Value = '  +4.39%  '
$ws.Range("D15").Value = '29.683.37'
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").Value = '3.81'
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("D17").Value = '64.32'
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("D18").Value = '241.32'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = '8.18'
$ws.Range("E19").Value = '  +7.84%  '
$ws.Range("D20").Value = '0.0₃0704'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '4.05'
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("E23").Value = '  +3.44%  '
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  +1.83%  '
$ws.Range("D25").Value = '156.51'
$ws.Range("E25").Value = '  +1.28%  '
$ws.Range("D26").Value = '15.51'
$ws.Range("E26").Value = '  +1.54%  '
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").Value = '6.50'
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("D30").Value = '0.0480'
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("E33").Value = '  +2.25%  '
$ws.Range("D34").Value = '1.426.53'
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  +4.44%  '
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("D37").Value = '2.84'
$ws.Range("E37").Value = '  +1.43%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  +2.74%  '
$ws.Range("D40").Value = '0.552'
$ws.Range("E40").Value = '  +3.24%  '
$ws.Range("D41").Value = '55.15'
$ws.Range("E41").Value = '  +1.98%  '
$ws.Range("D42").Value = '0.0494'
$ws.Range("E42").Value = '  +4.88%  '
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("D44").Value = '0.821'
$ws.Range("E44").Value = '  +3.43%  '
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").Value = '67.77'
$ws.Range("E46").Value = '  +3.28%  '
$ws.Range("D47").Value = '0.998'
$ws.Range("E47").Value = '  +19.47%  '
$ws.Range("E48").Value = '  +3.18%  '
$ws.Range("D49").Value = '1.743.16'
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("B50").Value = 'mCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D50").Value = '2.12'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '86.77'
$ws.Range("E51").Value = '  +0.33%  '